$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "transaction.delete.project.setProjectSection" right
# after the existing "transaction.delete.project.setProject" row (row 51),
# matching the format (style) of the surrounding data rows.
$ws.Rows(52).Insert()
$ws.Range("B55:C55").Copy()
$ws.Range("B52:C52").PasteSpecial(-4122)
$ws.Range("B52").Value = "transaction.delete.project.setProjectSection"
$ws.Range("C52").Value = "Menghapus Data Seksi Proyek"

# Insert a new row for "transaction.undelete.project.setProjectSection" right
# after the existing "transaction.undelete.project.setProject" row (now row
# 135, after the previous insert shifted everything down by one).
$ws.Rows(136).Insert()
$ws.Range("B135:C135").Copy()
$ws.Range("B136:C136").PasteSpecial(-4122)
$ws.Range("B136").Value = "transaction.undelete.project.setProjectSection"
$ws.Range("C136").Value = $ws.Range("C135").Value2

# Reflect the author's final cursor position/selection from the diff.
$ws.Range("B137").Select()

Write-Output "done"
